# Generate Report for Handoff
#
# The "b.md" row on every sheet moves from "Handed back: in sync with en-US"
# to "Ready for handoff", and the zh-cn / de-de detail sheets get the new
# handoff file name + handoff datetime for that same row.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("B3").Value = "Ready for handoff"

$hl = $ws.Range("C3").Hyperlinks.Item(1)
$hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"

$ws.Range("D3").Value = "2016-02-26 06:26:40"

# --- de-de sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("B3").Value = "Ready for handoff"

$hl = $ws.Range("C3").Hyperlinks.Item(1)
$hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"

$ws.Range("D3").Value = "2016-02-26 06:26:55"
